# "contingencies with rene fine"
# Populate the diagnostic sheet with a disconnected-elements count row/column
# header pair and apply the matching header formatting (bold, thin box
# border, centered/top-aligned) to B1 and A2. B2 holds the label text
# "disconnected_elements" (written through the shared-strings table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the header style once on B1 ...
$headerCell = $ws.Range("B1")
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop
$headerCell.Borders.LineStyle = 1         # xlContinuous (thin box border)

# ... then clone that exact formatting onto A2 via a formats-only paste so
# both header cells resolve to the same cell-style record.
$headerCell.Copy()
$ws.Range("A2").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false
